# FlixelRL-637: "毒の杖が強すぎる" (poison wand too strong)
# Rebalance consumable buy prices, equipment HP-ring comments, and item
# appearance (drop) rates.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("item_consumable")
$ws2 = $wb.Worksheets.Item("item_equipment")
$ws3 = $wb.Worksheets.Item("item_appear")

# ------------------------------------------------------------------
# item_consumable (sheet1): column J = buy price (K = J*0.35 sell, auto-recalcs)
# ------------------------------------------------------------------
$ws1.Range("J20").Value = 50
$ws1.Range("J21").Value = 100
$ws1.Range("J22").Value = 250
$ws1.Range("J23").Value = 500
$ws1.Range("J24").Value = 750
$ws1.Range("J25").Value = 1000
$ws1.Range("J26").Value = 1200
$ws1.Range("J27").Value = 1400
$ws1.Range("J28").Value = 1600
$ws1.Range("J29").Value = 1800
$ws1.Range("J30").Value = 2000
$ws1.Range("J31").Value = 2200
$ws1.Range("J32").Value = 2400
$ws1.Range("J33").Value = 2600
$ws1.Range("J34").Value = 2800

$ws1.Range("J67").Value = 3000
$ws1.Range("J69").Value = 2000
$ws1.Range("J74").Value = 1500

# ------------------------------------------------------------------
# item_equipment (sheet2): HP-ring comments bumped up one tier each
# (RING15 +10->+20, RING16 +20->+30, RING17 +30->+40, RING18 +40->+50)
# ------------------------------------------------------------------
$ws2.Range("M57").Value = "最大HPが20上昇します"
$ws2.Range("M58").Value = "最大HPが30上昇します"
$ws2.Range("M59").Value = "最大HPが40上昇します"
$ws2.Range("M60").Value = "最大HPが50上昇します"

# ------------------------------------------------------------------
# item_appear (sheet3): column B = start, column E = ratio
# ------------------------------------------------------------------
$ws3.Range("E42").Value = 30
$ws3.Range("E44").Value = 50
$ws3.Range("E46").Value = 50

$ws3.Range("B49").Value = 10
$ws3.Range("E49").Value = 50

$ws3.Range("E50").Value = 50

$ws3.Range("B51").Value = 5
$ws3.Range("B52").Value = 15

$ws3.Range("E58").Value = 0
$ws3.Range("E59").Value = 0
$ws3.Range("E60").Value = 10
$ws3.Range("E62").Value = 50

$ws3.Range("B66").Value = 10
$ws3.Range("E67").Value = 30
$ws3.Range("B69").Value = 15

$ws3.Range("B73").Value = 40
$ws3.Range("B74").Value = 40
$ws3.Range("B75").Value = 40
